$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 199.75
$ws.Range("I12").Value = 199
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 199
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = -29
$ws.Range("N12").Value = -540
$ws.Range("H70").Value = 15243634
$ws.Range("I70").Value = 33534194
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 100602582
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -100602312
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 15243634
$ws.Range("I73").Value = 33534194
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 100602582
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -100601646
$ws.Range("N73").Value = -6372
$ws.Range("H80").Value = 5927.25
$ws.Range("I80").Value = 215.77777
$ws.Range("K80").Value = 647.33331
$ws.Range("M80").Value = 350.66669
$ws.Range("H83").Value = 5927.25
$ws.Range("I83").Value = 215.77777
$ws.Range("K83").Value = 1941.99993
$ws.Range("M83").Value = 3050.00007
$ws.Range("H118").Value = 1289.1765
$ws.Range("I118").Value = 656.61536
$ws.Range("J118").Value = 3345
$ws.Range("K118").Value = 1969.84608
$ws.Range("L118").Value = 10035
$ws.Range("M118").Value = -312.84608
$ws.Range("N118").Value = -13349
$ws.Range("H121").Value = 2157
$ws.Range("J121").Value = 2996.6667
$ws.Range("L121").Value = 8990.000100000001
$ws.Range("N121").Value = -12484.0001
$ws.Range("H125").Value = 3314.6191
$ws.Range("J125").Value = 3072.6667
$ws.Range("L125").Value = 27654.0003
$ws.Range("N125").Value = -32574.0003
$ws.Range("H127").Value = 2033.4706
$ws.Range("I127").Value = 798
$ws.Range("J127").Value = 2548.25
$ws.Range("K127").Value = 2394
$ws.Range("L127").Value = 7644.75
$ws.Range("M127").Value = 2566
$ws.Range("N127").Value = -17564.75
$ws.Range("H138").Value = 1760.5555
$ws.Range("I138").Value = 1115.5135
$ws.Range("J138").Value = 4743.875
$ws.Range("K138").Value = 3346.5405
$ws.Range("L138").Value = 14231.625
$ws.Range("M138").Value = 1793.4595
$ws.Range("N138").Value = -24511.625

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 6000
$ws.Range("J46").Value = 6000
$ws.Range("L46").Value = 6000
$ws.Range("N46").Value = -6638
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").ClearContents()
$ws.Range("H122").Value = 2180.913
$ws.Range("I122").Value = 1524.0667
$ws.Range("J122").Value = 3412.5
$ws.Range("K122").Value = 4572.2001
$ws.Range("L122").Value = 10237.5
$ws.Range("M122").Value = -2122.2001
$ws.Range("N122").Value = -15137.5
$ws.Range("H132").Value = 8198.182000000001
$ws.Range("I132").Value = 11528.833
$ws.Range("J132").Value = 4201.4
$ws.Range("K132").Value = 34586.499
$ws.Range("L132").Value = 12604.2
$ws.Range("M132").Value = -32056.499
$ws.Range("N132").Value = -17664.2

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 47577.184
$ws.Range("I86").Value = 2276.2942
$ws.Range("J86").Value = 201600.2
$ws.Range("K86").Value = 2276.2942
$ws.Range("L86").Value = 201600.2
$ws.Range("M86").Value = -1153.2942
$ws.Range("N86").Value = -203846.2
$ws.Range("H89").Value = 47577.184
$ws.Range("I89").Value = 2276.2942
$ws.Range("J89").Value = 201600.2
$ws.Range("K89").Value = 11381.471
$ws.Range("L89").Value = 1008001
$ws.Range("M89").Value = -5765.471
$ws.Range("N89").Value = -1019233
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1101.4546
$ws.Range("I16").Value = 1071.5714
$ws.Range("J16").Value = 1153.75
$ws.Range("K16").Value = 1071.5714
$ws.Range("L16").Value = 1153.75
$ws.Range("M16").Value = -784.5714
$ws.Range("N16").Value = -1727.75
$ws.Range("H113").Value = 1101.4546
$ws.Range("I113").Value = 1071.5714
$ws.Range("J113").Value = 1153.75
$ws.Range("K113").Value = 1071.5714
$ws.Range("L113").Value = 1153.75
$ws.Range("M113").Value = 1098.4286
$ws.Range("N113").Value = -5493.75

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 177.75
$ws.Range("I40").Value = 60.285713
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 241.142852
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -172.142852
$ws.Range("N40").Value = -4138
$ws.Range("H49").Value = 5002.6665
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 5002.6665
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 15007.9995
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -15319.9995
$ws.Range("H103").Value = 1674.65
$ws.Range("J103").Value = 2697.0908
$ws.Range("L103").Value = 8091.2724
$ws.Range("N103").Value = -9849.2724

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 35000
$ws.Range("J32").Value = 35000
$ws.Range("L32").Value = 35000
$ws.Range("N32").Value = -35592
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H102").Value = 3877.7778
$ws.Range("I102").Value = 3980
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 3980
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -2358
$ws.Range("N102").Value = -6994
$ws.Range("I107").Value = 2090.6667
$ws.Range("J107").Value = 150
$ws.Range("K107").Value = 2090.6667
$ws.Range("L107").Value = 150
$ws.Range("M107").Value = -170.6667000000002
$ws.Range("N107").Value = -3990
$ws.Range("H114").Value = 38361
$ws.Range("J114").Value = 38361
$ws.Range("L114").Value = 38361
$ws.Range("N114").Value = -47039
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H120").Value = 20317
$ws.Range("J120").Value = 20317
$ws.Range("L120").Value = 20317
$ws.Range("N120").Value = -29993

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 21801.092
$ws.Range("I61").Value = 28501.5
$ws.Range("J61").Value = 3933.3333
$ws.Range("K61").Value = 28501.5
$ws.Range("L61").Value = 3933.3333
$ws.Range("M61").Value = -28299.5
$ws.Range("N61").Value = -4337.3333
$ws.Range("H113").Value = 21801.092
$ws.Range("I113").Value = 28501.5
$ws.Range("J113").Value = 3933.3333
$ws.Range("K113").Value = 28501.5
$ws.Range("L113").Value = 3933.3333
$ws.Range("M113").Value = -26331.5
$ws.Range("N113").Value = -8273.3333

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 262450
$ws.Range("J63").Value = 262450
$ws.Range("L63").Value = 262450
$ws.Range("N63").Value = -263698
$ws.Range("H66").Value = 262450
$ws.Range("J66").Value = 262450
$ws.Range("L66").Value = 787350
$ws.Range("N66").Value = -793590
$ws.Range("H81").Value = 220759.8
$ws.Range("I81").Value = 220759.8
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 441519.6
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -440458.6
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 220759.8
$ws.Range("I84").Value = 220759.8
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 2207598
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -2202294
$ws.Range("N84").ClearContents()
$ws.Range("H93").Value = 25389
$ws.Range("J93").Value = 25389
$ws.Range("L93").Value = 25389
$ws.Range("N93").Value = -30381
$ws.Range("H100").Value = 8741.68
$ws.Range("I100").Value = 14694.429
$ws.Range("J100").Value = 1165.4546
$ws.Range("K100").Value = 29388.858
$ws.Range("L100").Value = 2330.9092
$ws.Range("M100").Value = -28847.858
$ws.Range("N100").Value = -3412.9092
$ws.Range("H112").Value = 275000
$ws.Range("J112").Value = 275000
$ws.Range("L112").Value = 275000
$ws.Range("N112").Value = -277954
$ws.Range("H113").Value = 677.9375
$ws.Range("I113").Value = 594.2
$ws.Range("J113").Value = 817.5
$ws.Range("K113").Value = 1782.6
$ws.Range("L113").Value = 2452.5
$ws.Range("M113").Value = 387.3999999999999
$ws.Range("N113").Value = -6792.5
